# Update for Grey independent and other adjustments
# New polling figures roll in: the latest poll's numbers move into row 6
# (keeping the NSW/B and QLD/E columns, which are unaffected this round),
# and the previously "latest"/"second" rows cascade down to "second"/"third".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calc")

# Row 6 - "Latest Morgan ->" : only D (WA), F (National), G (SA+WA+etc) change
$ws.Range("D6").Value = 58
$ws.Range("F6").Value = 63.5
$ws.Range("G6").Value = 53

# Row 7 - "Second Morgan ->" : picks up what used to be row 6's figures
$ws.Range("B7").Value = 57
$ws.Range("C7").Value = 55
$ws.Range("D7").Value = 60.5
$ws.Range("E7").Value = 50.5
$ws.Range("F7").Value = 59
$ws.Range("G7").Value = 56

# Row 8 - "Third Morgan ->" : picks up what used to be row 7's figures
$ws.Range("B8").Value = 55.5
$ws.Range("C8").Value = 53
$ws.Range("D8").Value = 60
$ws.Range("E8").Value = 49
$ws.Range("F8").Value = 57
$ws.Range("G8").Value = 63.5

# Update the active selection to reflect where the editor was working
$ws.Activate()
$ws.Range("I18").Select()
